$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Special Items")
$ws.Range("A10").Value = "guitar"
